$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-21 Thursday" "2024-03-22 Friday"

Replace-Text "582÷3=194, 0" "885÷4=221, 1"
Replace-Text "598÷4=149, 2" "669÷3=223, 0"
Replace-Text "299÷7=42, 5" "246÷6=41, 0"
Replace-Text "593÷9=65, 8" "796÷7=113, 5"
Replace-Text "289÷4=72, 1" "719÷4=179, 3"

Replace-Text "944÷3=314, 2" "580÷4=145, 0"
Replace-Text "114÷7=16, 2" "978÷6=163, 0"
Replace-Text "780÷5=156, 0" "291÷9=32, 3"
Replace-Text "115÷9=12, 7" "776÷3=258, 2"
Replace-Text "696÷6=116, 0" "821÷8=102, 5"

Replace-Text "167÷5=33, 2" "697÷6=116, 1"
Replace-Text "572÷6=95, 2" "724÷3=241, 1"
Replace-Text "813÷6=135, 3" "816÷9=90, 6"
Replace-Text "670÷8=83, 6" "671÷2=335, 1"
Replace-Text "150÷6=25, 0" "851÷6=141, 5"

Replace-Text "396÷2=198, 0" "180÷9=20, 0"
Replace-Text "741÷8=92, 5" "585÷9=65, 0"
Replace-Text "207÷6=34, 3" "442÷7=63, 1"
Replace-Text "850÷9=94, 4" "671÷8=83, 7"
Replace-Text "917÷6=152, 5" "134÷8=16, 6"

Replace-Text "411÷2=205, 1" "748÷7=106, 6"
Replace-Text "491÷9=54, 5" "956÷4=239, 0"
Replace-Text "635÷7=90, 5" "544÷7=77, 5"
Replace-Text "629÷8=78, 5" "638÷5=127, 3"
Replace-Text "522÷5=104, 2" "715÷2=357, 1"

Write-Output "Done applying replacements"
